$d = $word.ActiveDocument

# 1. Bump the letter's date by a day: 3/15/24 -> 3/16/24
$found1 = $d.Content.Find.Execute("3/15/24", $true, $false, $false, $false, $false,
                                   $true, 1, $false, "3/16/24", 2)
if (-not $found1) {
    throw "Could not find the date '3/15/24' to replace."
}

# 2. Richard's wording suggestion: soften/rephrase the sentence following the
#    footnote reference about the paper's problems.
$old = ", published in your journal, as I feel both you and your readership should be aware of the problems before basing any decisions on the claims made. In particular, the paper asserts two claims:"
$new = ", published in your journal. It is important both you and your readership are aware of its shortcomings before basing any decisions on the claims made. In particular, the paper asserts two claims:"
$found2 = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                                   $true, 1, $false, $new, 2)
if (-not $found2) {
    throw "Could not find the sentence about the paper's shortcomings to replace."
}
